$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New Pathways data (replaces rows 2-14 and adds new row 15).
# Columns: A=Name, B=AID, C=Species, D=Location, E=Intake Date
$data = @(
    @("Marble",   "57806699", "Cat",               "Cat Adoption Room G",     "4/28/2025"),
    @("Goldie",   "58622481", "Cat",               "Dog Adoptions C",         "6/2/2025"),
    @("Trivento", "58448688", "Cat",               "Cat Adoption Condo Rooms","5/6/2025"),
    @("BEATRICE", "58442198", "Cat",               "Foster Home",             "5/6/2025"),
    @("Luna",     "58744079", "Mammal",             "Farm",                    "6/17/2025"),
    @("DASHWOOD", "58470388", "Cat",               "Cat Adoption Condo Rooms","5/9/2025"),
    @("DYA",      "58811207", "Dog",               "Dog Adoptions C",         "6/30/2025"),
    @("Glow",     "58834563", "Dog",               "Dog Holding E",           "7/2/2025"),
    @("Katniss",  "58834486", "Dog",               "Dog Holding E",           "7/2/2025"),
    @("Dior",     "58834490", "Dog",               "Foster Home",             "7/2/25"),
    @("Remy",     "58834525", "Dog",               "Foster Home",             "7/2/25"),
    @("Hilda",    "58917913", "Reptile/Amphibian", "Small Animals & Exotics", "7/16/25"),
    @("SMORES",   "58710884", "Cat",               "Foster Home",             "6/14/25"),
    @("Remy",     "58959672", "Dog",               "Dog Adoptions D",         "7/31/25")
)

$firstRow = 2
$lastRow = $firstRow + $data.Count - 1

# Format the AID and Intake Date columns as text first so values that look
# like numbers/dates (e.g. "58622481" or "4/28/2025") are stored as plain
# text, matching the source data (shared strings), not coerced numbers/dates.
$ws.Range("B$firstRow`:B$lastRow").NumberFormat = "@"
$ws.Range("E$firstRow`:E$lastRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $firstRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}

# Drop the temporary text-format styling so the cells end up with no
# explicit style (matching the original, un-styled data cells).
$ws.Range("B$firstRow`:B$lastRow").Style = "Normal"
$ws.Range("E$firstRow`:E$lastRow").Style = "Normal"
